# Burn Down Chart.xlsx edit script
# 1. Insert a new column before column A (shifts Timeline/Ideal/Actual columns to B/C/D)
# 2. Populate the new column A with a "Date" header and a series of dates (formatted as short date)
# 3. Resize the new column A
# 4. Re-establish the shared formula in (new) column C so the grouping survives the insert
# 5. Repoint the chart series formulas at their new columns (C/D instead of B/C)
# 6. Nudge the chart's anchor so it still starts one column to the right of the data table

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert new column A, shifting existing data to B:D ---
$ws.Columns.Item(1).Insert()

# --- 2. New column A: header + date values ---
$ws.Range("A1").Value = "Date"

# Apply the date number format to A2 first, then propagate it to A3:A21 via copy
# (keeps every date cell pointing at the same style index instead of creating
# a duplicate style per cell).
$ws.Range("A2").Value = 42982
$ws.Range("A2").NumberFormat = "mm-dd-yy"
$ws.Range("A2").Copy($ws.Range("A3:A21"))

$dates = @(42982,42983,42984,42985,42986,42989,42990,42991,42992,42993,42996,42997,42998,42999,43000,43003,43004,43005,43006,43007)
for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $ws.Range("A" + $row).Value = $dates[$i]
}

# --- 3. Width of the new column A ---
$ws.Columns.Item(1).ColumnWidth = 14.5

# --- 4. Rebuild the shared formula block in column C (Insert() flattens it) ---
$ws.Range("C4:C21").Formula = '=ROUND($C$2-($C$2/20*B3)-2,0)'

# --- 5. Fix chart series references (now one column further right) ---
$co = $ws.ChartObjects().Item(1)
$s1 = $co.Chart.SeriesCollection().Item(1)
$s2 = $co.Chart.SeriesCollection().Item(2)
$s1.Formula = '=SERIES(''Sprint 1''!$C$1,,''Sprint 1''!$C$2:$C$21,1)'
$s2.Formula = '=SERIES(''Sprint 1''!$D$1,,''Sprint 1''!$D$2:$D$21,2)'

# --- 6. Shift chart anchor right by the width of the new column A ---
$co.Left = $co.Left + $ws.Columns.Item(1).Width
